$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the old "_GoBack" bookmark that sat inside the
#    "System administrators can now create and edit groups." paragraph.
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------------
# 2. The previously-empty paragraph right after that sentence becomes a new
#    highlighted TODO item.
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs(4)
$p4.Range.Text = "TODO: Add error handling to all update methods"

$todoAdd = $d.Content.Duplicate
$todoAdd.Find.Execute("TODO: Add error handling to all update methods")
$todoAdd.HighlightColorIndex = 7

# ---------------------------------------------------------------------------
# 3. Completed TODO items get struck through (whole paragraph, including the
#    paragraph mark, so Word records the formatting on both <w:pPr>/<w:rPr>
#    and on every run).
# ---------------------------------------------------------------------------
$d.Paragraphs(11).Range.Font.StrikeThrough = 1   # TODO: Hide the button ...
$d.Paragraphs(12).Range.Font.StrikeThrough = 1   # TODO: Handle the parsing ...

# ---------------------------------------------------------------------------
# 4. The still-open TODO item gets highlighted instead (only the visible
#    text, not the paragraph mark).
# ---------------------------------------------------------------------------
$todoExtend = $d.Content.Duplicate
$todoExtend.Find.Execute("TODO: Extend the get list of events controller to return past events")
$todoExtend.HighlightColorIndex = 7

# ---------------------------------------------------------------------------
# 5. "There are now some missing tests (frowning face)" and the following
#    "When the user first registers ..." paragraph both get highlighted,
#    including their paragraph marks. Rebuild them precisely via InsertXML
#    so the emoji run's existing <w:rPr> keeps its <mc:AlternateContent> and
#    simply gains a <w:highlight>.
# ---------------------------------------------------------------------------
$xmlMissingTests = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex"><w:pPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">There are now some missing tests </w:t></w:r><w:r><w:rPr><mc:AlternateContent><mc:Choice Requires="w16se"/><mc:Fallback><w:rFonts w:ascii="Segoe UI Emoji" w:eastAsia="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/></mc:Fallback></mc:AlternateContent><w:highlight w:val="yellow"/></w:rPr><mc:AlternateContent><mc:Choice Requires="w16se"><w16se:symEx w16se:font="Segoe UI Emoji" w16se:char="2639"/></mc:Choice><mc:Fallback><w:t>&#x2639;</w:t></mc:Fallback></mc:AlternateContent></w:r></w:p>'
$d.Paragraphs(22).Range.InsertXML($xmlMissingTests)

$xmlFirstRegisters = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:tab/><w:t>When the user first registers, we need a test to ensure that a message is sent to slack</w:t></w:r></w:p>'
$d.Paragraphs(23).Range.InsertXML($xmlFirstRegisters)

# ---------------------------------------------------------------------------
# 6. The slack-webhook TODO gets highlighted too, but only the run itself
#    (no paragraph-mark formatting).
# ---------------------------------------------------------------------------
$slackWebhook = $d.Content.Duplicate
$slackWebhook.Find.Execute("We need a test for the slack webhook")
$slackWebhook.HighlightColorIndex = 7

# ---------------------------------------------------------------------------
# 7. Completed work items further down get struck through as well.
# ---------------------------------------------------------------------------
$d.Paragraphs(39).Range.Font.StrikeThrough = 1   # Add ApprovedBy / RejectedBy ...
$d.Paragraphs(40).Range.Font.StrikeThrough = 1   # Add screen to allow the user ...

# ---------------------------------------------------------------------------
# 8. The new "latest edit" location gets the _GoBack bookmark (collapsed,
#    right before the run) plus a highlight on that run.
# ---------------------------------------------------------------------------
$generic = $d.Content.Duplicate
$generic.Find.Execute("It would be nice to make this more generic so that in the future the group owner can define their own questions.")
$genericStart = $generic.Start
$d.Range($genericStart, $genericStart).Bookmarks.Add("_GoBack")

$generic2 = $d.Content.Duplicate
$generic2.Find.Execute("It would be nice to make this more generic so that in the future the group owner can define their own questions.")
$generic2.HighlightColorIndex = 7
